# Fruta / hortaliza, semanal
# Shuffle the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values across the
# existing data rows (2-22) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    3  = @{ D = 44432; M = 30;  N = 1300; O = 1300; P = 1300; S = 1300 }
    4  = @{ D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
    5  = @{ D = 44763; M = 50;  N = 2300; O = 2300; P = 2300; S = 2300 }
    6  = @{ D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    7  = @{ D = 44476; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    8  = @{ D = 44812; M = 50;  N = 2500; O = 2500; P = 2500; S = 2500 }
    9  = @{ D = 44418; M = 40;  N = 1200; O = 1200; P = 1200; S = 1200 }
    10 = @{ D = 44830; M = 50;  N = 2500; O = 2500; P = 2500; S = 2500 }
    11 = @{ D = 44811; M = 60;  N = 2500; O = 2500; P = 2500; S = 2500 }
    12 = @{ D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 }
    13 = @{ D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    14 = @{ D = 44762; M = 50;  N = 2300; O = 2300; P = 2300; S = 2300 }
    15 = @{ D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 }
    16 = @{ D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 }
    17 = @{ D = 44749; M = 120; N = 2300; O = 2300; P = 2300; S = 2300 }
    18 = @{ D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
    19 = @{ D = 44748; M = 300; N = 2300; O = 2300; P = 2300; S = 2300 }
    20 = @{ D = 44760; M = 80;  N = 2300; O = 2300; P = 2300; S = 2300 }
    21 = @{ D = 44753; M = 160; N = 2300; O = 2300; P = 2300; S = 2300 }
    22 = @{ D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 }
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values.D
    $ws.Range("M$row").Value = $values.M
    $ws.Range("N$row").Value = $values.N
    $ws.Range("O$row").Value = $values.O
    $ws.Range("P$row").Value = $values.P
    $ws.Range("S$row").Value = $values.S
}
